$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and the Kaspa/ThetaToken row swap)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.146.62"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.984.70"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  +11.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.75"
$ws.Range("E6").Value = "  +8.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.681"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.750"
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.47"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000319"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.91"
$ws.Range("E13").Value = "  +3.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.627.20"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.988.68"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.27"
$ws.Range("E16").Value = "  +9.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.30"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.830.92"
$ws.Range("E20").Value = "  +2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.63"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("E22").Value = "  +13.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.88"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("E24").Value = "  -4.32%  "
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  +15.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.41"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.27"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.75"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.69"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "47.98"
$ws.Range("E34").Value = "  -5.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "664.82"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.92"
$ws.Range("E36").Value = "  +9.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0906"
$ws.Range("E37").Value = "  +11.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.436"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.34"
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.145"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0488"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.60"
$ws.Range("E45").Value = "  +6.74%  "
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.45"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.875.05"
$ws.Range("E49").Value = "  +9.49%  "
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("E51").Value = "  +4.59%  "

Write-Output "Applied cryptos update"
